$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (FAPs -> Ptn/Sdc3 -> ECs)
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 3.270036666666666
$ws.Range("H2").Value = 9.81011
$ws.Range("I2").Value = 0.359406393324744
$ws.Range("J2").Value = 0.3594063933247441
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 23.59622066666667
$ws.Range("N2").Value = 70.788662
$ws.Range("O2").Value = 0.6996728317814862
$ws.Range("P2").Value = 0.6996728317814862
$ws.Range("Q2").Value = 77.16050677475778
$ws.Range("R2").Value = 694.44456097282
$ws.Range("S2").Value = 0.2514668889778943
$ws.Range("T2").Value = 0.2514668889778943

# Row 3 (FAPs -> Ptn/Sdc3 -> FAPs)
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 3.270036666666666
$ws.Range("H3").Value = 9.81011
$ws.Range("I3").Value = 0.359406393324744
$ws.Range("J3").Value = 0.3594063933247441
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 7.778025666666667
$ws.Range("N3").Value = 23.334077
$ws.Range("O3").Value = 0.2306332577891816
$ws.Range("P3").Value = 0.2306332577891816
$ws.Range("Q3").Value = 25.43442912427444
$ws.Range("R3").Value = 228.90986211847
$ws.Range("S3").Value = 0.08289106736274569
$ws.Range("T3").Value = 0.08289106736274571

# Row 4 (FAPs -> Ptn/Sdc3 -> sCs)
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 3.270036666666666
$ws.Range("H4").Value = 9.81011
$ws.Range("I4").Value = 0.359406393324744
$ws.Range("J4").Value = 0.3594063933247441
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 2.350402666666667
$ws.Range("N4").Value = 7.051208000000001
$ws.Range("O4").Value = 0.06969391042933218
$ws.Range("P4").Value = 0.06969391042933218
$ws.Range("Q4").Value = 7.685902901431112
$ws.Range("R4").Value = 69.17312611288001
$ws.Range("S4").Value = 0.02504843698410404
$ws.Range("T4").Value = 0.02504843698410404

# Row 5 (sCs -> Ptn/Sdc3 -> ECs)
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 5.828401
$ws.Range("H5").Value = 17.485203
$ws.Range("I5").Value = 0.6405936066752559
$ws.Range("J5").Value = 0.640593606675256
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 23.59622066666667
$ws.Range("N5").Value = 70.788662
$ws.Range("O5").Value = 0.6996728317814862
$ws.Range("P5").Value = 0.6996728317814862
$ws.Range("Q5").Value = 137.5282361298207
$ws.Range("R5").Value = 1237.754125168386
$ws.Range("S5").Value = 0.4482059428035918
$ws.Range("T5").Value = 0.4482059428035919

# Row 6 (sCs -> Ptn/Sdc3 -> FAPs)
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 5.828401
$ws.Range("H6").Value = 17.485203
$ws.Range("I6").Value = 0.6405936066752559
$ws.Range("J6").Value = 0.640593606675256
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 7.778025666666667
$ws.Range("N6").Value = 23.334077
$ws.Range("O6").Value = 0.2306332577891816
$ws.Range("P6").Value = 0.2306332577891816
$ws.Range("Q6").Value = 45.33345257362567
$ws.Range("R6").Value = 408.001073162631
$ws.Range("S6").Value = 0.1477421904264359
$ws.Range("T6").Value = 0.1477421904264359

# Row 7 (sCs -> Ptn/Sdc3 -> sCs)
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 5.828401
$ws.Range("H7").Value = 17.485203
$ws.Range("I7").Value = 0.6405936066752559
$ws.Range("J7").Value = 0.640593606675256
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 2.350402666666667
$ws.Range("N7").Value = 7.051208000000001
$ws.Range("O7").Value = 0.06969391042933218
$ws.Range("P7").Value = 0.06969391042933218
$ws.Range("Q7").Value = 13.69908925280267
$ws.Range("R7").Value = 123.291803275224
$ws.Range("S7").Value = 0.04464547344522813
$ws.Range("T7").Value = 0.04464547344522814
